$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data rows (header row 1 is left untouched, preserving its
# formatting/style) so the shared-strings table rebuilds cleanly for the
# rewritten data, in the same string order the source tool produced.
$ws.Range("A2:T13").Clear()

# Sending cluster column first (establishes new cluster-name string order)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(17, 1).Value = "sCs"

# Target cluster column
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(17, 4).Value = "sCs"

# Ligand symbol column
$ws.Cells.Item(2, 2).Value = "Ntn1"
$ws.Cells.Item(3, 2).Value = "Ntn1"
$ws.Cells.Item(4, 2).Value = "Ntn1"
$ws.Cells.Item(5, 2).Value = "Ntn1"
$ws.Cells.Item(6, 2).Value = "Ntn1"
$ws.Cells.Item(7, 2).Value = "Ntn1"
$ws.Cells.Item(8, 2).Value = "Ntn1"
$ws.Cells.Item(9, 2).Value = "Ntn1"
$ws.Cells.Item(10, 2).Value = "Ntn1"
$ws.Cells.Item(11, 2).Value = "Ntn1"
$ws.Cells.Item(12, 2).Value = "Ntn1"
$ws.Cells.Item(13, 2).Value = "Ntn1"
$ws.Cells.Item(14, 2).Value = "Ntn1"
$ws.Cells.Item(15, 2).Value = "Ntn1"
$ws.Cells.Item(16, 2).Value = "Ntn1"
$ws.Cells.Item(17, 2).Value = "Ntn1"

# Receptor symbol column
$ws.Cells.Item(2, 3).Value = "Unc5a"
$ws.Cells.Item(3, 3).Value = "Unc5a"
$ws.Cells.Item(4, 3).Value = "Unc5a"
$ws.Cells.Item(5, 3).Value = "Unc5a"
$ws.Cells.Item(6, 3).Value = "Unc5a"
$ws.Cells.Item(7, 3).Value = "Unc5a"
$ws.Cells.Item(8, 3).Value = "Unc5a"
$ws.Cells.Item(9, 3).Value = "Unc5a"
$ws.Cells.Item(10, 3).Value = "Unc5a"
$ws.Cells.Item(11, 3).Value = "Unc5a"
$ws.Cells.Item(12, 3).Value = "Unc5a"
$ws.Cells.Item(13, 3).Value = "Unc5a"
$ws.Cells.Item(14, 3).Value = "Unc5a"
$ws.Cells.Item(15, 3).Value = "Unc5a"
$ws.Cells.Item(16, 3).Value = "Unc5a"
$ws.Cells.Item(17, 3).Value = "Unc5a"

# Remaining numeric columns
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.990837
$ws.Cells.Item(3, 7).Value = 1.990837
$ws.Cells.Item(4, 7).Value = 1.990837
$ws.Cells.Item(5, 7).Value = 1.990837
$ws.Cells.Item(6, 7).Value = 11.42765333333333
$ws.Cells.Item(7, 7).Value = 11.42765333333333
$ws.Cells.Item(8, 7).Value = 11.42765333333333
$ws.Cells.Item(9, 7).Value = 11.42765333333333
$ws.Cells.Item(10, 7).Value = 0.2266433333333333
$ws.Cells.Item(11, 7).Value = 0.2266433333333333
$ws.Cells.Item(12, 7).Value = 0.2266433333333333
$ws.Cells.Item(13, 7).Value = 0.2266433333333333
$ws.Cells.Item(14, 7).Value = 4.085149333333333
$ws.Cells.Item(15, 7).Value = 4.085149333333333
$ws.Cells.Item(16, 7).Value = 4.085149333333333
$ws.Cells.Item(17, 7).Value = 4.085149333333333
$ws.Cells.Item(2, 8).Value = 5.972511000000001
$ws.Cells.Item(3, 8).Value = 5.972511000000001
$ws.Cells.Item(4, 8).Value = 5.972511000000001
$ws.Cells.Item(5, 8).Value = 5.972511000000001
$ws.Cells.Item(6, 8).Value = 34.28296
$ws.Cells.Item(7, 8).Value = 34.28296
$ws.Cells.Item(8, 8).Value = 34.28296
$ws.Cells.Item(9, 8).Value = 34.28296
$ws.Cells.Item(10, 8).Value = 0.67993
$ws.Cells.Item(11, 8).Value = 0.67993
$ws.Cells.Item(12, 8).Value = 0.67993
$ws.Cells.Item(13, 8).Value = 0.67993
$ws.Cells.Item(14, 8).Value = 12.255448
$ws.Cells.Item(15, 8).Value = 12.255448
$ws.Cells.Item(16, 8).Value = 12.255448
$ws.Cells.Item(17, 8).Value = 12.255448
$ws.Cells.Item(2, 9).Value = 0.1122845585713437
$ws.Cells.Item(3, 9).Value = 0.1122845585713437
$ws.Cells.Item(4, 9).Value = 0.1122845585713437
$ws.Cells.Item(5, 9).Value = 0.1122845585713437
$ws.Cells.Item(6, 9).Value = 0.6445274073365515
$ws.Cells.Item(7, 9).Value = 0.6445274073365515
$ws.Cells.Item(8, 9).Value = 0.6445274073365515
$ws.Cells.Item(9, 9).Value = 0.6445274073365515
$ws.Cells.Item(10, 9).Value = 0.01278283789002879
$ws.Cells.Item(11, 9).Value = 0.01278283789002879
$ws.Cells.Item(12, 9).Value = 0.01278283789002879
$ws.Cells.Item(13, 9).Value = 0.01278283789002879
$ws.Cells.Item(14, 9).Value = 0.230405196202076
$ws.Cells.Item(15, 9).Value = 0.230405196202076
$ws.Cells.Item(16, 9).Value = 0.230405196202076
$ws.Cells.Item(17, 9).Value = 0.230405196202076
$ws.Cells.Item(2, 10).Value = 0.1122845585713437
$ws.Cells.Item(3, 10).Value = 0.1122845585713437
$ws.Cells.Item(4, 10).Value = 0.1122845585713437
$ws.Cells.Item(5, 10).Value = 0.1122845585713437
$ws.Cells.Item(6, 10).Value = 0.6445274073365515
$ws.Cells.Item(7, 10).Value = 0.6445274073365515
$ws.Cells.Item(8, 10).Value = 0.6445274073365515
$ws.Cells.Item(9, 10).Value = 0.6445274073365515
$ws.Cells.Item(10, 10).Value = 0.01278283789002879
$ws.Cells.Item(11, 10).Value = 0.01278283789002879
$ws.Cells.Item(12, 10).Value = 0.01278283789002879
$ws.Cells.Item(13, 10).Value = 0.01278283789002879
$ws.Cells.Item(14, 10).Value = 0.230405196202076
$ws.Cells.Item(15, 10).Value = 0.230405196202076
$ws.Cells.Item(16, 10).Value = 0.230405196202076
$ws.Cells.Item(17, 10).Value = 0.230405196202076
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.837459
$ws.Cells.Item(3, 13).Value = 1.864510666666667
$ws.Cells.Item(4, 13).Value = 0.4771786666666666
$ws.Cells.Item(5, 13).Value = 1.320881333333333
$ws.Cells.Item(6, 13).Value = 0.837459
$ws.Cells.Item(7, 13).Value = 1.864510666666667
$ws.Cells.Item(8, 13).Value = 0.4771786666666666
$ws.Cells.Item(9, 13).Value = 1.320881333333333
$ws.Cells.Item(10, 13).Value = 0.837459
$ws.Cells.Item(11, 13).Value = 1.864510666666667
$ws.Cells.Item(12, 13).Value = 0.4771786666666666
$ws.Cells.Item(13, 13).Value = 1.320881333333333
$ws.Cells.Item(14, 13).Value = 0.837459
$ws.Cells.Item(15, 13).Value = 1.864510666666667
$ws.Cells.Item(16, 13).Value = 0.4771786666666666
$ws.Cells.Item(17, 13).Value = 1.320881333333333
$ws.Cells.Item(2, 14).Value = 2.512377
$ws.Cells.Item(3, 14).Value = 5.593532
$ws.Cells.Item(4, 14).Value = 1.431536
$ws.Cells.Item(5, 14).Value = 3.962644
$ws.Cells.Item(6, 14).Value = 2.512377
$ws.Cells.Item(7, 14).Value = 5.593532
$ws.Cells.Item(8, 14).Value = 1.431536
$ws.Cells.Item(9, 14).Value = 3.962644
$ws.Cells.Item(10, 14).Value = 2.512377
$ws.Cells.Item(11, 14).Value = 5.593532
$ws.Cells.Item(12, 14).Value = 1.431536
$ws.Cells.Item(13, 14).Value = 3.962644
$ws.Cells.Item(14, 14).Value = 2.512377
$ws.Cells.Item(15, 14).Value = 5.593532
$ws.Cells.Item(16, 14).Value = 1.431536
$ws.Cells.Item(17, 14).Value = 3.962644
$ws.Cells.Item(2, 15).Value = 0.1861007731134217
$ws.Cells.Item(3, 15).Value = 0.4143329721752205
$ws.Cells.Item(4, 15).Value = 0.1060390046317472
$ws.Cells.Item(5, 15).Value = 0.2935272500796106
$ws.Cells.Item(6, 15).Value = 0.1861007731134217
$ws.Cells.Item(7, 15).Value = 0.4143329721752205
$ws.Cells.Item(8, 15).Value = 0.1060390046317472
$ws.Cells.Item(9, 15).Value = 0.2935272500796106
$ws.Cells.Item(10, 15).Value = 0.1861007731134217
$ws.Cells.Item(11, 15).Value = 0.4143329721752205
$ws.Cells.Item(12, 15).Value = 0.1060390046317472
$ws.Cells.Item(13, 15).Value = 0.2935272500796106
$ws.Cells.Item(14, 15).Value = 0.1861007731134217
$ws.Cells.Item(15, 15).Value = 0.4143329721752205
$ws.Cells.Item(16, 15).Value = 0.1060390046317472
$ws.Cells.Item(17, 15).Value = 0.2935272500796106
$ws.Cells.Item(2, 16).Value = 0.1861007731134217
$ws.Cells.Item(3, 16).Value = 0.4143329721752205
$ws.Cells.Item(4, 16).Value = 0.1060390046317472
$ws.Cells.Item(5, 16).Value = 0.2935272500796106
$ws.Cells.Item(6, 16).Value = 0.1861007731134217
$ws.Cells.Item(7, 16).Value = 0.4143329721752205
$ws.Cells.Item(8, 16).Value = 0.1060390046317472
$ws.Cells.Item(9, 16).Value = 0.2935272500796106
$ws.Cells.Item(10, 16).Value = 0.1861007731134217
$ws.Cells.Item(11, 16).Value = 0.4143329721752205
$ws.Cells.Item(12, 16).Value = 0.1060390046317472
$ws.Cells.Item(13, 16).Value = 0.2935272500796106
$ws.Cells.Item(14, 16).Value = 0.1861007731134217
$ws.Cells.Item(15, 16).Value = 0.4143329721752205
$ws.Cells.Item(16, 16).Value = 0.1060390046317472
$ws.Cells.Item(17, 16).Value = 0.2935272500796106
$ws.Cells.Item(2, 17).Value = 1.667244363183
$ws.Cells.Item(3, 17).Value = 3.711936822094667
$ws.Cells.Item(4, 17).Value = 0.9499849452106667
$ws.Cells.Item(5, 17).Value = 2.629659431009334
$ws.Cells.Item(6, 17).Value = 9.57019113288
$ws.Cells.Item(7, 17).Value = 21.30698153496889
$ws.Cells.Item(8, 17).Value = 5.453032380728889
$ws.Cells.Item(9, 17).Value = 15.09457397180445
$ws.Cells.Item(10, 17).Value = 0.18980449929
$ws.Cells.Item(11, 17).Value = 0.4225789125288889
$ws.Cells.Item(12, 17).Value = 0.1081493636088889
$ws.Cells.Item(13, 17).Value = 0.2993689483244444
$ws.Cells.Item(14, 17).Value = 3.421145075544
$ws.Cells.Item(15, 17).Value = 7.616804506926221
$ws.Cells.Item(16, 17).Value = 1.949346112014222
$ws.Cells.Item(17, 17).Value = 5.395997498279111
$ws.Cells.Item(2, 18).Value = 15.005199268647
$ws.Cells.Item(3, 18).Value = 33.407431398852
$ws.Cells.Item(4, 18).Value = 8.549864506896
$ws.Cells.Item(5, 18).Value = 23.666934879084
$ws.Cells.Item(6, 18).Value = 86.13172019592
$ws.Cells.Item(7, 18).Value = 191.76283381472
$ws.Cells.Item(8, 18).Value = 49.07729142656
$ws.Cells.Item(9, 18).Value = 135.85116574624
$ws.Cells.Item(10, 18).Value = 1.70824049361
$ws.Cells.Item(11, 18).Value = 3.80321021276
$ws.Cells.Item(12, 18).Value = 0.97334427248
$ws.Cells.Item(13, 18).Value = 2.69432053492
$ws.Cells.Item(14, 18).Value = 30.790305679896
$ws.Cells.Item(15, 18).Value = 68.55124056233599
$ws.Cells.Item(16, 18).Value = 17.544115008128
$ws.Cells.Item(17, 18).Value = 48.563977484512
$ws.Cells.Item(2, 19).Value = 0.02089624315882635
$ws.Cells.Item(3, 19).Value = 0.04652319488224748
$ws.Cells.Item(4, 19).Value = 0.01190654282642041
$ws.Cells.Item(5, 19).Value = 0.0329585777038495
$ws.Cells.Item(6, 19).Value = 0.1199470487981215
$ws.Cells.Item(7, 19).Value = 0.2670489563301424
$ws.Cells.Item(8, 19).Value = 0.06834504473184862
$ws.Cells.Item(9, 19).Value = 0.189186357476439
$ws.Cells.Item(10, 19).Value = 0.002378896013917898
$ws.Cells.Item(11, 19).Value = 0.005296351215809653
$ws.Cells.Item(12, 19).Value = 0.001355479406227637
$ws.Cells.Item(13, 19).Value = 0.003752111254073603
$ws.Cells.Item(14, 19).Value = 0.04287858514255596
$ws.Cells.Item(15, 19).Value = 0.09546446974702097
$ws.Cells.Item(16, 19).Value = 0.02443193766725057
$ws.Cells.Item(17, 19).Value = 0.06763020364524851
$ws.Cells.Item(2, 20).Value = 0.02089624315882635
$ws.Cells.Item(3, 20).Value = 0.04652319488224748
$ws.Cells.Item(4, 20).Value = 0.01190654282642041
$ws.Cells.Item(5, 20).Value = 0.0329585777038495
$ws.Cells.Item(6, 20).Value = 0.1199470487981215
$ws.Cells.Item(7, 20).Value = 0.2670489563301424
$ws.Cells.Item(8, 20).Value = 0.06834504473184862
$ws.Cells.Item(9, 20).Value = 0.189186357476439
$ws.Cells.Item(10, 20).Value = 0.002378896013917898
$ws.Cells.Item(11, 20).Value = 0.005296351215809653
$ws.Cells.Item(12, 20).Value = 0.001355479406227637
$ws.Cells.Item(13, 20).Value = 0.003752111254073603
$ws.Cells.Item(14, 20).Value = 0.04287858514255596
$ws.Cells.Item(15, 20).Value = 0.09546446974702098
$ws.Cells.Item(16, 20).Value = 0.02443193766725057
$ws.Cells.Item(17, 20).Value = 0.06763020364524852
